# The commit simplifies word/styles.xml's <w:docDefaults> down to only the
# handful of character/paragraph properties that are NOT already Word's
# built-in defaults (Arial 11pt/en language; 1.15 "auto" line spacing).
# Everything else that was spelled out explicitly (bold/italic/smallCaps/
# strike/color/underline/shading/vertAlign off; keepNext/keepLines off,
# widowControl on, no borders/shading, zero spacing-before/after, zero
# indents, no contextual spacing, left justification) was redundant -
# those are exactly Word's own defaults, so re-asserting them through the
# object model below is a content/formatting no-op, matching the diff's
# intent of leaving the *effective* formatting of the document unchanged.

$d = $word.ActiveDocument
$normal = $d.Styles("Normal")

# --- rPrDefault (character formatting defaults) ---
$font = $normal.Font
$font.Name       = "Arial"
$font.NameAscii  = "Arial"
$font.NameOther  = "Arial"
$font.NameFarEast = "Arial"
$font.Size       = 11
$font.SizeBi     = 11
$font.Bold           = $false
$font.Italic         = $false
$font.SmallCaps      = $false
$font.StrikeThrough  = $false
$font.Color          = -16777216
$font.Underline      = 0
$font.Superscript    = $false
$font.Subscript      = $false
$font.Shading.Texture               = 0
$font.Shading.ForegroundPatternColor = -16777216
$font.Shading.BackgroundPatternColor = -16777216

# --- pPrDefault (paragraph formatting defaults) ---
$pf = $normal.ParagraphFormat
$pf.KeepWithNext   = $false
$pf.KeepTogether   = $false
$pf.WidowControl   = $true
$pf.SpaceBefore    = 0
$pf.SpaceAfter     = 0
$pf.LineSpacingRule = 5
$pf.LineSpacing     = 13.8
$pf.LeftIndent      = 0
$pf.RightIndent     = 0
$pf.FirstLineIndent = 0
$pf.Alignment       = 0
$normal.NoSpaceBetweenParagraphsOfSameStyle = $false

$border = $pf.Borders
$border.Item(-1).LineStyle = 0
$border.Item(1).LineStyle = 0
$border.Item(2).LineStyle = 0
$border.Item(3).LineStyle = 0
$border.Item(4).LineStyle = 0
$border.Item(5).LineStyle = 0

$pf.Shading.Texture               = 0
$pf.Shading.ForegroundPatternColor = -16777216
$pf.Shading.BackgroundPatternColor = -16777216

Write-Output "docDefaults normalized"
